$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 values ---
$ws.Range("A2").Value = "001/DR IFRAN"
$ws.Range("C2").Value = "DDDD"
$ws.Range("D2").Value = "ALI EXPRESSE"
$ws.Range("G2").Value = 12000
$ws.Range("H2").Value = 8000
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = 1200
$ws.Range("K2").Value = 800
$ws.Range("L2").Value = 4000
$ws.Range("M2").Value = 10800

# --- Add new row 3 ---
$ws.Range("A3").Value = "001/DR IFRAN"
$ws.Range("B3").Value = "Direction régionale"
$ws.Range("C3").Value = "BB132345"
$ws.Range("D3").Value = "KHALID TAGHMAOUI"
$ws.Range("E3").Value = "ds"
$ws.Range("F3").Value = "mensuelle"
$ws.Range("G3").Value = 9000
$ws.Range("H3").Value = 6000
$ws.Range("I3").Value = 10
$ws.Range("J3").Value = 900
$ws.Range("K3").Value = 600
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = 8100

# --- Add new row 4 ---
$ws.Range("A4").Value = "001/DR IFRAN"
$ws.Range("B4").Value = "Direction régionale"
$ws.Range("C4").Value = "KS123456"
$ws.Range("D4").Value = "Youssef You"
$ws.Range("E4").Value = "ds"
$ws.Range("F4").Value = "mensuelle"
$ws.Range("G4").Value = 9000
$ws.Range("H4").Value = 6000
$ws.Range("I4").Value = 10
$ws.Range("J4").Value = 900
$ws.Range("K4").Value = 600
$ws.Range("L4").Value = 3000
$ws.Range("M4").Value = 8100
